# Change share separator from ':' to '.' in the Model sheet's shareCodes
# column (B7:B9), e.g. "ZAR:ALSI" -> "ZAR.ALSI".
$wb = $excel.ActiveWorkbook

$wsCallSimple = $wb.Worksheets.Item("Script Call Simple")
$wsIncentive  = $wb.Worksheets.Item("Script Incentive Option")
$wsModel      = $wb.Worksheets.Item("Model")

$wsModel.Range("B7").Value = "ZAR.ALSI"
$wsModel.Range("B8").Value = "ZAR.AAA"
$wsModel.Range("B9").Value = "ZAR.BBB"

# Update the selections / active sheet to match the state the workbook was
# saved in: the user had moved off "Script Call Simple" (was tabSelected)
# and ended up with "Model" active, with new cell selections on each sheet.
$wsCallSimple.Range("D24").Select() | Out-Null
$wsIncentive.Range("C4").Select() | Out-Null

$wsModel.Activate() | Out-Null
$wsModel.Range("B10").Select() | Out-Null
